$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph near the top of the document.
#    (Paragraph 2: "Meta description" (bold) + ": Read our review of
#    Football Cash Collect, a 5x3 slot game with 30 paylines, four
#    jackpots, and free spin feature. Play Football Cash Collect for free
#    today!")
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2. The final paragraph (the AI image-prompt paragraph, in italics) gets
#    split into two paragraphs:
#      a) a new bold paragraph carrying the page title
#         "Play Football Cash Collect for Free - Review & Ratings"
#      b) the original (italic) paragraph, but its text becomes the former
#         meta-description text, minus the "Meta description: " label.
# ---------------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)
$lastRange = $lastPara.Range

$newHeading = "Play Football Cash Collect for Free - Review & Ratings"

# Insert the new heading text as a plain, unformatted run right before the
# existing (italic) text, then split the two apart into separate
# paragraphs by inserting a paragraph break right after the new text.
$lastRange.InsertBefore($newHeading)

$splitAt = $lastRange.Start + $newHeading.Length
$splitPoint = $d.Range($splitAt, $splitAt)
$splitPoint.InsertParagraphAfter()

# The heading text now lives in its own paragraph - make it bold (only the
# text run, not the trailing paragraph mark).
$headingPara = $d.Paragraphs.Item($paraCount)
$headingRange = $headingPara.Range
$headingTextRange = $d.Range($headingRange.Start, $headingRange.End - 1)
$headingTextRange.Bold = 1

# The original image-prompt text now lives in the paragraph right after the
# new heading; replace it with the former meta-description copy (without
# the "Meta description: " label), keeping its existing italic run.
$descPara = $d.Paragraphs.Item($paraCount + 1)
$descRange = $descPara.Range
$descTextRange = $d.Range($descRange.Start, $descRange.End - 1)
$descTextRange.Text = "Read our review of Football Cash Collect, a 5x3 slot game with 30 paylines, four jackpots, and free spin feature. Play Football Cash Collect for free today!"
